# "added create a car using excel"
# Clean up the "Model Year" header in C1: the original cell held rich text
# made of two runs ("Model Year" + an italic trailing non-breaking space).
# Re-assigning a plain string collapses it back down to a single, uniformly
# formatted run using the cell's existing style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 3).Value = "Model Year"

# Leave the selection on the header cell that was just touched.
$ws.Range("C1").Select()
